$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the "datetimeFigureOut" date placeholder text that lives on
#    the slide master and on every slide layout: 2020/10/19 -> 2020/10/30
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {}
        if ($isDatePh -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -match "2020/10/19") {
                $tr.Text = "2020/10/30"
            } else {
                $tr.Text = "2020/10/30"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders($master.Shapes)

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders($layouts.Item($li).Shapes)
}

# ---------------------------------------------------------------------
# 2) Rename the "Length" column header to "Offset" in the small 4-column
#    (Min / AvgInc / Length / BitRequired) tables used on slides 17-19.
# ---------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
                for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                    $cell = $tbl.Cell($r, $c)
                    $ctr = $cell.Shape.TextFrame.TextRange
                    if ($ctr.Text -eq "Length") {
                        $ctr.Text = "Offset"
                    }
                }
            }
        }
    }
}
